# Auto-generated edit script: refresh the charging-station idle-time report.
# Overwrites the data rows (now 37 data rows, rows 2-38) and the header
# label order in row 1, matching the re-pulled report for 2025-09-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: columns C/D/E relabelled (same text, new column order) ---
$ws.Cells.Item(1,3).Value = "最近一次充电结束时间"
$ws.Cells.Item(1,4).Value = "截止一直未充电时间"
$ws.Cells.Item(1,5).Value = "截止一直未充电时长(小时)"

# --- Data rows ---
$ws.Cells.Item(2,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(2,2).Value = "402号直流"
$ws.Cells.Item(2,3).Value = 45915.503680555557
$ws.Cells.Item(2,4).Value = 45928.332673611112
$ws.Cells.Item(2,5).Value = 307.89583333331393
$ws.Cells.Item(3,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(3,2).Value = "201号直流"
$ws.Cells.Item(3,3).Value = 45925.246076388888
$ws.Cells.Item(3,4).Value = 45928.332673611112
$ws.Cells.Item(3,5).Value = 74.078333333367482
$ws.Cells.Item(4,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(4,2).Value = "101号直流"
$ws.Cells.Item(4,3).Value = 45926.043692129628
$ws.Cells.Item(4,4).Value = 45928.332673611112
$ws.Cells.Item(4,5).Value = 54.935555555624887
$ws.Cells.Item(5,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(5,2).Value = "603号直流"
$ws.Cells.Item(5,3).Value = 45926.099652777775
$ws.Cells.Item(5,4).Value = 45928.332673611112
$ws.Cells.Item(5,5).Value = 53.592500000086147
$ws.Cells.Item(6,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(6,2).Value = "602号直流"
$ws.Cells.Item(6,3).Value = 45926.242071759261
$ws.Cells.Item(6,4).Value = 45928.332673611112
$ws.Cells.Item(6,5).Value = 50.17444444441935
$ws.Cells.Item(7,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(7,2).Value = "008B号直流"
$ws.Cells.Item(7,3).Value = 45926.525636574072
$ws.Cells.Item(7,4).Value = 45928.332673611112
$ws.Cells.Item(7,5).Value = 43.368888888973743
$ws.Cells.Item(8,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(8,2).Value = "001A号直流"
$ws.Cells.Item(8,3).Value = 45926.556597222225
$ws.Cells.Item(8,4).Value = 45928.332673611112
$ws.Cells.Item(8,5).Value = 42.625833333295304
$ws.Cells.Item(9,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(9,2).Value = "903号直流"
$ws.Cells.Item(9,3).Value = 45926.560879629629
$ws.Cells.Item(9,4).Value = 45928.332673611112
$ws.Cells.Item(9,5).Value = 42.523055555589963
$ws.Cells.Item(10,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(10,2).Value = "801号直流"
$ws.Cells.Item(10,3).Value = 45927.038657407407
$ws.Cells.Item(10,4).Value = 45928.332673611112
$ws.Cells.Item(10,5).Value = 31.056388888915535
$ws.Cells.Item(11,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(11,2).Value = "402号直流"
$ws.Cells.Item(11,3).Value = 45927.050925925927
$ws.Cells.Item(11,4).Value = 45928.332673611112
$ws.Cells.Item(11,5).Value = 30.761944444442634
$ws.Cells.Item(12,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(12,2).Value = "406号直流"
$ws.Cells.Item(12,3).Value = 45927.051412037035
$ws.Cells.Item(12,4).Value = 45928.332673611112
$ws.Cells.Item(12,5).Value = 30.750277777842712
$ws.Cells.Item(13,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(13,2).Value = "905号直流"
$ws.Cells.Item(13,3).Value = 45927.051574074074
$ws.Cells.Item(13,4).Value = 45928.332673611112
$ws.Cells.Item(13,5).Value = 30.746388888917863
$ws.Cells.Item(14,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(14,2).Value = "401号直流"
$ws.Cells.Item(14,3).Value = 45927.157442129632
$ws.Cells.Item(14,4).Value = 45928.332673611112
$ws.Cells.Item(14,5).Value = 28.205555555527098
$ws.Cells.Item(15,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(15,2).Value = "A04号直流"
$ws.Cells.Item(15,3).Value = 45927.251550925925
$ws.Cells.Item(15,4).Value = 45928.332673611112
$ws.Cells.Item(15,5).Value = 25.946944444498513
$ws.Cells.Item(16,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(16,2).Value = "003B号直流"
$ws.Cells.Item(16,3).Value = 45927.302870370368
$ws.Cells.Item(16,4).Value = 45928.332673611112
$ws.Cells.Item(16,5).Value = 24.715277777868323
$ws.Cells.Item(17,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(17,2).Value = "103号直流"
$ws.Cells.Item(17,3).Value = 45927.358912037038
$ws.Cells.Item(17,4).Value = 45928.332673611112
$ws.Cells.Item(17,5).Value = 23.370277777779847
$ws.Cells.Item(18,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18,2).Value = "502号直流"
$ws.Cells.Item(18,3).Value = 45927.408541666664
$ws.Cells.Item(18,4).Value = 45928.332673611112
$ws.Cells.Item(18,5).Value = 22.179166666755918
$ws.Cells.Item(19,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19,2).Value = "B01号直流"
$ws.Cells.Item(19,3).Value = 45927.413043981483
$ws.Cells.Item(19,4).Value = 45928.332673611112
$ws.Cells.Item(19,5).Value = 22.071111111086793
$ws.Cells.Item(20,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20,2).Value = "803号直流"
$ws.Cells.Item(20,3).Value = 45927.41479166667
$ws.Cells.Item(20,4).Value = 45928.332673611112
$ws.Cells.Item(20,5).Value = 22.02916666661622
$ws.Cells.Item(21,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21,2).Value = "904号直流"
$ws.Cells.Item(21,3).Value = 45927.445706018516
$ws.Cells.Item(21,4).Value = 45928.332673611112
$ws.Cells.Item(21,5).Value = 21.287222222308628
$ws.Cells.Item(22,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22,2).Value = "701号直流"
$ws.Cells.Item(22,3).Value = 45927.457337962966
$ws.Cells.Item(22,4).Value = 45928.332673611112
$ws.Cells.Item(22,5).Value = 21.008055555517785
$ws.Cells.Item(23,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(23,2).Value = "B02号直流"
$ws.Cells.Item(23,3).Value = 45927.508564814816
$ws.Cells.Item(23,4).Value = 45928.332673611112
$ws.Cells.Item(23,5).Value = 19.77861111110542
$ws.Cells.Item(24,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24,2).Value = "204号直流"
$ws.Cells.Item(24,3).Value = 45927.523240740738
$ws.Cells.Item(24,4).Value = 45928.332673611112
$ws.Cells.Item(24,5).Value = 19.426388888969086
$ws.Cells.Item(25,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(25,2).Value = "102号直流"
$ws.Cells.Item(25,3).Value = 45927.524837962963
$ws.Cells.Item(25,4).Value = 45928.332673611112
$ws.Cells.Item(25,5).Value = 19.38805555558065
$ws.Cells.Item(26,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(26,2).Value = "504号直流"
$ws.Cells.Item(26,3).Value = 45927.534305555557
$ws.Cells.Item(26,4).Value = 45928.332673611112
$ws.Cells.Item(26,5).Value = 19.160833333327901
$ws.Cells.Item(27,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(27,2).Value = "905号直流"
$ws.Cells.Item(27,3).Value = 45927.543576388889
$ws.Cells.Item(27,4).Value = 45928.332673611112
$ws.Cells.Item(27,5).Value = 18.938333333353512
$ws.Cells.Item(28,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(28,2).Value = "902号直流"
$ws.Cells.Item(28,3).Value = 45927.545995370368
$ws.Cells.Item(28,4).Value = 45928.332673611112
$ws.Cells.Item(28,5).Value = 18.880277777847368
$ws.Cells.Item(29,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(29,2).Value = "901号直流"
$ws.Cells.Item(29,3).Value = 45927.563981481479
$ws.Cells.Item(29,4).Value = 45928.332673611112
$ws.Cells.Item(29,5).Value = 18.448611111205537
$ws.Cells.Item(30,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(30,2).Value = "306号直流"
$ws.Cells.Item(30,3).Value = 45927.569039351853
$ws.Cells.Item(30,4).Value = 45928.332673611112
$ws.Cells.Item(30,5).Value = 18.327222222229466
$ws.Cells.Item(31,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(31,2).Value = "202号直流"
$ws.Cells.Item(31,3).Value = 45927.586284722223
$ws.Cells.Item(31,4).Value = 45928.332673611112
$ws.Cells.Item(31,5).Value = 17.913333333330229
$ws.Cells.Item(32,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(32,2).Value = "106号直流"
$ws.Cells.Item(32,3).Value = 45927.606226851851
$ws.Cells.Item(32,4).Value = 45928.332673611112
$ws.Cells.Item(32,5).Value = 17.434722222271375
$ws.Cells.Item(33,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(33,2).Value = "705号直流"
$ws.Cells.Item(33,3).Value = 45927.65315972222
$ws.Cells.Item(33,4).Value = 45928.332673611112
$ws.Cells.Item(33,5).Value = 16.308333333407063
$ws.Cells.Item(34,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(34,2).Value = "902号直流"
$ws.Cells.Item(34,3).Value = 45927.677430555559
$ws.Cells.Item(34,4).Value = 45928.332673611112
$ws.Cells.Item(34,5).Value = 15.725833333272021
$ws.Cells.Item(35,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(35,2).Value = "103号直流"
$ws.Cells.Item(35,3).Value = 45927.746261574073
$ws.Cells.Item(35,4).Value = 45928.332673611112
$ws.Cells.Item(35,5).Value = 14.073888888931833
$ws.Cells.Item(36,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(36,2).Value = "401号直流"
$ws.Cells.Item(36,3).Value = 45927.758240740739
$ws.Cells.Item(36,4).Value = 45928.332673611112
$ws.Cells.Item(36,5).Value = 13.786388888955116
$ws.Cells.Item(37,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(37,2).Value = "002A号直流"
$ws.Cells.Item(37,3).Value = 45927.776076388887
$ws.Cells.Item(37,4).Value = 45928.332673611112
$ws.Cells.Item(37,5).Value = 13.358333333395422
$ws.Cells.Item(38,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(38,2).Value = "005A号直流"
$ws.Cells.Item(38,3).Value = 45927.828460648147
$ws.Cells.Item(38,4).Value = 45928.332673611112
$ws.Cells.Item(38,5).Value = 12.101111111172941

# --- Column widths: nudge to match the workbook re-saved after the refresh ---
# (Excel recalculates "best fit" widths off real glyph metrics on save; these
#  are the closest values reachable through the ColumnWidth property.)
$ws.Columns.Item(1).ColumnWidth = 41.4296875
$ws.Columns.Item(2).ColumnWidth = 10.572544642857142
$ws.Columns.Item(3).ColumnWidth = 20.858258928571427
$ws.Columns.Item(4).ColumnWidth = 20.858258928571427
$ws.Columns.Item(5).ColumnWidth = 24.143973214285715

# --- Selection cursor moved to G5 (matches the saved sheetView) ---
$ws.Range("G5").Select() | Out-Null
